$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.044746123766386
$ws.Cells.Item(2, 4).Value = 1.046124982441258
$ws.Cells.Item(2, 5).Value = 1.052897441569706
$ws.Cells.Item(2, 6).Value = 1.064543961853475
$ws.Cells.Item(2, 9).Value = 1.045669872579522
$ws.Cells.Item(2, 10).Value = 1.04980955987066
$ws.Cells.Item(2, 11).Value = 1.048891295850239
$ws.Cells.Item(2, 12).Value = 1.055644903836449
$ws.Cells.Item(2, 13).Value = 1.067259613238838
$ws.Cells.Item(2, 14).Value = 1.020432980523155
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.045739685798793
$ws.Cells.Item(3, 4).Value = 1.046872203366356
$ws.Cells.Item(3, 5).Value = 1.05381357995693
$ws.Cells.Item(3, 6).Value = 1.065617123146953
$ws.Cells.Item(3, 9).Value = 1.045978034621589
$ws.Cells.Item(3, 10).Value = 1.050450284518014
$ws.Cells.Item(3, 11).Value = 1.049450263464028
$ws.Cells.Item(3, 12).Value = 1.056373714702255
$ws.Cells.Item(3, 13).Value = 1.06814735425734
$ws.Cells.Item(3, 14).Value = 1.020650040175898
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.04638279349944
$ws.Cells.Item(4, 4).Value = 1.047355814031524
$ws.Cells.Item(4, 5).Value = 1.054406949553625
$ws.Cells.Item(4, 6).Value = 1.066312365911983
$ws.Cells.Item(4, 9).Value = 1.04617624167151
$ws.Cells.Item(4, 10).Value = 1.050864470575411
$ws.Cells.Item(4, 11).Value = 1.049811396221749
$ws.Cells.Item(4, 12).Value = 1.056845239182677
$ws.Cells.Item(4, 13).Value = 1.068722015647432
$ws.Cells.Item(4, 14).Value = 1.020790255574705
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.046653205014375
$ws.Cells.Item(5, 4).Value = 1.047559149055476
$ws.Cells.Item(5, 5).Value = 1.054656537009399
$ws.Cells.Item(5, 6).Value = 1.066604845329782
$ws.Cells.Item(5, 9).Value = 1.046259281601756
$ws.Cells.Item(5, 10).Value = 1.051038496709283
$ws.Cells.Item(5, 11).Value = 1.049963082554154
$ws.Cells.Item(5, 12).Value = 1.057043451908691
$ws.Cells.Item(5, 13).Value = 1.068963658566696
$ws.Cells.Item(5, 14).Value = 1.020849145179868
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.046698611160803
$ws.Cells.Item(6, 4).Value = 1.047593291345358
$ws.Cells.Item(6, 5).Value = 1.054698451727066
$ws.Cells.Item(6, 6).Value = 1.066653965581374
$ws.Cells.Item(6, 9).Value = 1.046273207571091
$ws.Cells.Item(6, 10).Value = 1.051067710717957
$ws.Cells.Item(6, 11).Value = 1.049988543488324
$ws.Cells.Item(6, 12).Value = 1.05707673173128
$ws.Cells.Item(6, 13).Value = 1.069004234698379
$ws.Cells.Item(6, 14).Value = 1.0208590296561
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.046386406559177
$ws.Cells.Item(7, 4).Value = 1.047358530906614
$ws.Cells.Item(7, 5).Value = 1.054410284024466
$ws.Cells.Item(7, 6).Value = 1.066316273253795
$ws.Cells.Item(7, 9).Value = 1.046177352380599
$ws.Cells.Item(7, 10).Value = 1.050866796305678
$ws.Cells.Item(7, 11).Value = 1.049813423588763
$ws.Cells.Item(7, 12).Value = 1.05684788777513
$ws.Cells.Item(7, 13).Value = 1.068725244275577
$ws.Cells.Item(7, 14).Value = 1.020791042684609
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.045081859664116
$ws.Cells.Item(8, 4).Value = 1.046377486168294
$ws.Cells.Item(8, 5).Value = 1.053206936976639
$ws.Cells.Item(8, 6).Value = 1.064906468291165
$ws.Cells.Item(8, 9).Value = 1.045774264468042
$ws.Cells.Item(8, 10).Value = 1.050026179408065
$ws.Cells.Item(8, 11).Value = 1.049080316465524
$ws.Cells.Item(8, 12).Value = 1.055891221934982
$ws.Cells.Item(8, 13).Value = 1.067559580714132
$ws.Cells.Item(8, 14).Value = 1.020506385693599
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.042784682627534
$ws.Cells.Item(9, 4).Value = 1.044649634180904
$ws.Cells.Item(9, 5).Value = 1.051090862589609
$ws.Cells.Item(9, 6).Value = 1.062428642464861
$ws.Cells.Item(9, 9).Value = 1.045054840769236
$ws.Cells.Item(9, 10).Value = 1.048541823598763
$ws.Cells.Item(9, 11).Value = 1.047784251530519
$ws.Cells.Item(9, 12).Value = 1.054204985323632
$ws.Cells.Item(9, 13).Value = 1.065507350184254
$ws.Cells.Item(9, 14).Value = 1.020002983066026
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.041254330897243
$ws.Cells.Item(10, 4).Value = 1.043498374906441
$ws.Cells.Item(10, 5).Value = 1.049683135836195
$ws.Cells.Item(10, 6).Value = 1.060781121439832
$ws.Cells.Item(10, 9).Value = 1.044569106774361
$ws.Cells.Item(10, 10).Value = 1.047550213899362
$ws.Cells.Item(10, 11).Value = 1.0469173948703
$ws.Cells.Item(10, 12).Value = 1.053080548119266
$ws.Cells.Item(10, 13).Value = 1.064140455350129
$ws.Cells.Item(10, 14).Value = 1.019666186786628
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.040591936634627
$ws.Cells.Item(11, 4).Value = 1.043000030260387
$ws.Cells.Item(11, 5).Value = 1.049074293027108
$ws.Cells.Item(11, 6).Value = 1.060068769096338
$ws.Cells.Item(11, 9).Value = 1.044357331595211
$ws.Cells.Item(11, 10).Value = 1.047120358705439
$ws.Cells.Item(11, 11).Value = 1.046541376056247
$ws.Cells.Item(11, 12).Value = 1.05259359522119
$ws.Cells.Item(11, 13).Value = 1.063548880571527
$ws.Cells.Item(11, 14).Value = 1.019520070206491
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.040345932933205
$ws.Cells.Item(12, 4).Value = 1.042814947480913
$ws.Cells.Item(12, 5).Value = 1.048848249387299
$ws.Cells.Item(12, 6).Value = 1.059804325776762
$ws.Cells.Item(12, 9).Value = 1.04427845142779
$ws.Cells.Item(12, 10).Value = 1.046960619438143
$ws.Cells.Item(12, 11).Value = 1.046401606649942
$ws.Cells.Item(12, 12).Value = 1.052412710103866
$ws.Cells.Item(12, 13).Value = 1.063329188894875
$ws.Cells.Item(12, 14).Value = 1.019465753941381
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.040398699792586
$ws.Cells.Item(13, 4).Value = 1.042854647238206
$ws.Cells.Item(13, 5).Value = 1.048896731631093
$ws.Cells.Item(13, 6).Value = 1.059861042695793
$ws.Cells.Item(13, 9).Value = 1.044295381326373
$ws.Cells.Item(13, 10).Value = 1.046994887299958
$ws.Cells.Item(13, 11).Value = 1.046431592146391
$ws.Cells.Item(13, 12).Value = 1.052451510980459
$ws.Cells.Item(13, 13).Value = 1.06337631143648
$ws.Cells.Item(13, 14).Value = 1.019477406864122
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.040571601094514
$ws.Cells.Item(14, 4).Value = 1.042984730753184
$ws.Cells.Item(14, 5).Value = 1.049055605988501
$ws.Cells.Item(14, 6).Value = 1.060046906942361
$ws.Cells.Item(14, 9).Value = 1.044350815769963
$ws.Cells.Item(14, 10).Value = 1.047107156072064
$ws.Cells.Item(14, 11).Value = 1.0465298246958
$ws.Cells.Item(14, 12).Value = 1.052578643384237
$ws.Cells.Item(14, 13).Value = 1.063530719857127
$ws.Cells.Item(14, 14).Value = 1.019515581263717
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.040678136470963
$ws.Cells.Item(15, 4).Value = 1.043064882786376
$ws.Cells.Item(15, 5).Value = 1.049153508014128
$ws.Cells.Item(15, 6).Value = 1.060161444734464
$ws.Cells.Item(15, 9).Value = 1.044384941956259
$ws.Cells.Item(15, 10).Value = 1.047176319052863
$ws.Cells.Item(15, 11).Value = 1.046590335871631
$ws.Cells.Item(15, 12).Value = 1.052656972659701
$ws.Cells.Item(15, 13).Value = 1.063625862030445
$ws.Cells.Item(15, 14).Value = 1.019539096204895
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.041298297281815
$ws.Cells.Item(16, 4).Value = 1.043531451804551
$ws.Cells.Item(16, 5).Value = 1.049723557791911
$ws.Cells.Item(16, 6).Value = 1.060828419773971
$ws.Cells.Item(16, 9).Value = 1.04458313108758
$ws.Cells.Item(16, 10).Value = 1.047578731863739
$ws.Cells.Item(16, 11).Value = 1.046942336050679
$ws.Cells.Item(16, 12).Value = 1.053112864271002
$ws.Cells.Item(16, 13).Value = 1.064179722577669
$ws.Cells.Item(16, 14).Value = 1.019675878151074
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.041687377051612
$ws.Cells.Item(17, 4).Value = 1.04382416117587
$ws.Cells.Item(17, 5).Value = 1.050081326112894
$ws.Cells.Item(17, 6).Value = 1.061247073110523
$ws.Cells.Item(17, 9).Value = 1.044707062112599
$ws.Cells.Item(17, 10).Value = 1.0478310261566
$ws.Cells.Item(17, 11).Value = 1.047162959005741
$ws.Cells.Item(17, 12).Value = 1.053398816214884
$ws.Cells.Item(17, 13).Value = 1.064527225263914
$ws.Cells.Item(17, 14).Value = 1.019761602638409
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.041914345588711
$ws.Cells.Item(18, 4).Value = 1.043994908809226
$ws.Cells.Item(18, 5).Value = 1.050290074988554
$ws.Cells.Item(18, 6).Value = 1.061491366357463
$ws.Cells.Item(18, 9).Value = 1.044779209141336
$ws.Cells.Item(18, 10).Value = 1.047978138646354
$ws.Cells.Item(18, 11).Value = 1.047291580527009
$ws.Cells.Item(18, 12).Value = 1.053565600920124
$ws.Cells.Item(18, 13).Value = 1.06472994658441
$ws.Cells.Item(18, 14).Value = 1.019811577098059
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.04199174020905
$ws.Cells.Item(19, 4).Value = 1.044053131893865
$ws.Cells.Item(19, 5).Value = 1.050361264586632
$ws.Cells.Item(19, 6).Value = 1.061574680959533
$ws.Cells.Item(19, 9).Value = 1.044803785685602
$ws.Cells.Item(19, 10).Value = 1.048028292302178
$ws.Cells.Item(19, 11).Value = 1.047335426251356
$ws.Cells.Item(19, 12).Value = 1.053622469122996
$ws.Cells.Item(19, 13).Value = 1.064799074199513
$ws.Cells.Item(19, 14).Value = 1.019828612479223
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.04164562989822
$ws.Cells.Item(20, 4).Value = 1.043792754651213
$ws.Cells.Item(20, 5).Value = 1.050042933850323
$ws.Cells.Item(20, 6).Value = 1.06120214523883
$ws.Cells.Item(20, 9).Value = 1.044693779948867
$ws.Cells.Item(20, 10).Value = 1.047803962173339
$ws.Cells.Item(20, 11).Value = 1.047139294878713
$ws.Cells.Item(20, 12).Value = 1.053368136917736
$ws.Cells.Item(20, 13).Value = 1.064489938523181
$ws.Cells.Item(20, 14).Value = 1.019752408022162
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.040520684891441
$ws.Cells.Item(21, 4).Value = 1.042946423714141
$ws.Cells.Item(21, 5).Value = 1.049008818471784
$ws.Cells.Item(21, 6).Value = 1.059992170241953
$ws.Cells.Item(21, 9).Value = 1.044334497708028
$ws.Cells.Item(21, 10).Value = 1.047074097698292
$ws.Cells.Item(21, 11).Value = 1.046500900394829
$ws.Cells.Item(21, 12).Value = 1.052541206307916
$ws.Cells.Item(21, 13).Value = 1.063485249164658
$ws.Cells.Item(21, 14).Value = 1.019504341013671
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.039813613200806
$ws.Cells.Item(22, 4).Value = 1.042414444669891
$ws.Cells.Item(22, 5).Value = 1.048359252484754
$ws.Cells.Item(22, 6).Value = 1.059232314113431
$ws.Cells.Item(22, 9).Value = 1.0441073448238
$ws.Cells.Item(22, 10).Value = 1.046614787117437
$ws.Cells.Item(22, 11).Value = 1.046098942313046
$ws.Cells.Item(22, 12).Value = 1.052021229647478
$ws.Cells.Item(22, 13).Value = 1.062853824800302
$ws.Cells.Item(22, 14).Value = 1.01934812820492
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.040188423790122
$ws.Cells.Item(23, 4).Value = 1.042696443009675
$ws.Cells.Item(23, 5).Value = 1.048703540382
$ws.Cells.Item(23, 6).Value = 1.059635042431807
$ws.Cells.Item(23, 9).Value = 1.044227882003149
$ws.Cells.Item(23, 10).Value = 1.046858315609521
$ws.Cells.Item(23, 11).Value = 1.046312082149319
$ws.Cells.Item(23, 12).Value = 1.052296883924819
$ws.Cells.Item(23, 13).Value = 1.063188529635336
$ws.Cells.Item(23, 14).Value = 1.019430962545756
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.041664493562865
$ws.Cells.Item(24, 4).Value = 1.043806945860877
$ws.Cells.Item(24, 5).Value = 1.050060281449044
$ws.Cells.Item(24, 6).Value = 1.061222445901636
$ws.Cells.Item(24, 9).Value = 1.044699782019004
$ws.Cells.Item(24, 10).Value = 1.047816191366089
$ws.Cells.Item(24, 11).Value = 1.047149987877213
$ws.Cells.Item(24, 12).Value = 1.053381999590768
$ws.Cells.Item(24, 13).Value = 1.064506786707061
$ws.Cells.Item(24, 14).Value = 1.019756562757467
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.043378365836823
$ws.Cells.Item(25, 4).Value = 1.045096216142378
$ws.Cells.Item(25, 5).Value = 1.051637395047453
$ws.Cells.Item(25, 6).Value = 1.063068453735587
$ws.Cells.Item(25, 9).Value = 1.045241908970701
$ws.Cells.Item(25, 10).Value = 1.048925926964464
$ws.Cells.Item(25, 11).Value = 1.048119813825468
$ws.Cells.Item(25, 12).Value = 1.054640969803992
$ws.Cells.Item(25, 13).Value = 1.066037681790218
$ws.Cells.Item(25, 14).Value = 1.020133336436214